$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the new L column values (2020 data)
$ws.Range("L4").Value = 2020
$ws.Range("L5").Value = 1.2
$ws.Range("L6").Value = 1.7
$ws.Range("L7").Value = 0.4
$ws.Range("L8").Value = 3.3
$ws.Range("L9").Value = 3.9
$ws.Range("L10").Value = 2.4
$ws.Range("L11").Value = 95.5
$ws.Range("L12").Value = 94.4
$ws.Range("L13").Value = 97.2

# Copy styles from column K so the new L cells match formatting of their row
$ws.Range("K4").Copy()
$ws.Range("L4").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("K5").Copy()
$ws.Range("L5").PasteSpecial(-4122)

$ws.Range("K6").Copy()
$ws.Range("L6").PasteSpecial(-4122)

$ws.Range("K7").Copy()
$ws.Range("L7").PasteSpecial(-4122)

$ws.Range("K8").Copy()
$ws.Range("L8").PasteSpecial(-4122)

$ws.Range("K9").Copy()
$ws.Range("L9").PasteSpecial(-4122)

$ws.Range("K10").Copy()
$ws.Range("L10").PasteSpecial(-4122)

$ws.Range("K11").Copy()
$ws.Range("L11").PasteSpecial(-4122)

$ws.Range("K12").Copy()
$ws.Range("L12").PasteSpecial(-4122)

$ws.Range("K13").Copy()
$ws.Range("L13").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Update the view: scroll so column C is leftmost, and select L4:L13 with active cell L4
$ws.Range("L4:L13").Select()
$excel.ActiveWindow.ScrollColumn = 3
